$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B498").Value = "Miscellaneous Options"
$ws.Range("C498").Value = "Opzioni varie"
